# Update odds values in row 3 of Sheet1 per the latest FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G3"  = 1.95
    "H3"  = 3.3
    "I3"  = 3.55
    "J3"  = 2.55
    "K3"  = 2.15
    "L3"  = 4
    "S3"  = 1.37
    "T3"  = 2.85
    "U3"  = 1.62
    "V3"  = 2.15
    "W3"  = 8.25
    "Z3"  = 18.5
    "AA3" = 14.5
    "AD3" = 6.6
    "AE3" = 12.5
    "AH3" = 12
    "AL3" = 30
    "AM3" = 32
    "AN3" = 4
    "AO3" = 10
    "AP3" = 17
    "AQ3" = 37
    "AR3" = 60
    "AS3" = 200
    "AT3" = 2.85
    "AU3" = 6.7
    "AW3" = 5.6
    "AY3" = 24
    "AZ3" = 100
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
